$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 897-940 (2 new rows appended at the end; a new
# "2026/02/27" day block is inserted, displacing the existing time/ranking (C/D)
# values by two rows while dates/weekdays (A/B) keep their own sequence).
$data = @(
    @(897, "2026/02/27", "金", 19, 201),
    @(898, "2026/02/27", "金", 22, 201),
    @(899, "2026/02/27", "金", 13, 88),
    @(900, "2026/02/27", "金", 16, 99),
    @(901, "2026/12/29", "火", 19, 81),
    @(902, "2026/12/29", "火", 23, 85),
    @(903, "2026/12/29", "火", 2, 89),
    @(904, "2026/12/29", "火", 5, 95),
    @(905, "2026/12/29", "火", 8, 91),
    @(906, "2026/12/29", "火", 13, 92),
    @(907, "2026/12/30", "水", 16, 99),
    @(908, "2026/12/30", "水", 22, 108),
    @(909, "2026/12/30", "水", 2, 114),
    @(910, "2026/12/30", "水", 6, 120),
    @(911, "2026/12/30", "水", 9, 120),
    @(912, "2026/12/30", "水", 12, 201),
    @(913, "2026/12/31", "木", 14, 130),
    @(914, "2026/12/31", "木", 22, 120),
    @(915, "2026/12/31", "木", 2, 129),
    @(916, "2026/12/31", "木", 5, 119),
    @(917, "2026/12/31", "木", 13, 133),
    @(918, "2027/01/01", "金", 16, 109),
    @(919, "2027/01/01", "金", 19, 120),
    @(920, "2027/01/01", "金", 1, 105),
    @(921, "2027/01/01", "金", 5, 109),
    @(922, "2027/01/01", "金", 8, 110),
    @(923, "2027/01/01", "金", 13, 132),
    @(924, "2027/01/01", "金", 16, 145),
    @(925, "2027/01/02", "土", 19, 157),
    @(926, "2027/01/02", "土", 22, 165),
    @(927, "2027/01/02", "土", 1, 174),
    @(928, "2027/01/02", "土", 4, 192),
    @(929, "2027/01/02", "土", 7, 189),
    @(930, "2027/01/02", "土", 13, 201),
    @(931, "2027/01/02", "土", 16, 201),
    @(932, "2027/01/03", "日", 19, 201),
    @(933, "2027/01/03", "日", 22, 194),
    @(934, "2027/01/03", "日", 2, 164),
    @(935, "2027/01/03", "日", 5, 166),
    @(936, "2027/01/03", "日", 7, 168),
    @(937, "2027/01/04", "月", 13, 173),
    @(938, "2027/01/04", "月", 22, 127),
    @(939, "2027/01/05", "火", 1, 118),
    @(940, "2027/01/05", "火", 7, 127)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateVal = $row[1]
    $wdVal = $row[2]
    $timeVal = $row[3]
    $rankVal = $row[4]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $dateVal
    $aCell.ClearFormats()

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $wdVal
    $bCell.ClearFormats()

    $ws.Cells.Item($r, 3).Value = $timeVal
    $ws.Cells.Item($r, 4).Value = $rankVal
}
